$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage and preserving the
# cells original style (avoids Excel auto-converting numeric-looking
# strings like "5.60" or "0.0000146" into real numbers, which would drop
# formatting such as trailing zeros).
function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "63.722.46"
Set-TextValue $ws.Range("E2") "  +0.26%  "
Set-TextValue $ws.Range("D3") "2.653.95"
Set-TextValue $ws.Range("E3") "  +0.33%  "
Set-TextValue $ws.Range("E4") "  +0.07%  "
Set-TextValue $ws.Range("D5") "604.01"
Set-TextValue $ws.Range("E5") "  +2.21%  "
Set-TextValue $ws.Range("D6") "147.85"
Set-TextValue $ws.Range("E6") "  +2.33%  "
Set-TextValue $ws.Range("E7") "  +0.02%  "
Set-TextValue $ws.Range("E8") "  -0.08%  "
Set-TextValue $ws.Range("D9") "0.109"
Set-TextValue $ws.Range("E9") "  +1.46%  "
Set-TextValue $ws.Range("D10") "5.60"
Set-TextValue $ws.Range("E10") "  +0.24%  "
Set-TextValue $ws.Range("D11") "0.370"
Set-TextValue $ws.Range("E11") "  +4.76%  "
Set-TextValue $ws.Range("E12") "  +0.18%  "
Set-TextValue $ws.Range("D13") "27.57"
Set-TextValue $ws.Range("E13") "  +0.43%  "
Set-TextValue $ws.Range("D14") "3.128.38"
Set-TextValue $ws.Range("E14") "  +0.39%  "
Set-TextValue $ws.Range("D15") "63.584.24"
Set-TextValue $ws.Range("E15") "  +0.23%  "
Set-TextValue $ws.Range("D16") "0.0000146"
Set-TextValue $ws.Range("E16") "  +0.48%  "
Set-TextValue $ws.Range("D17") "2.645.49"
Set-TextValue $ws.Range("E17") "  +0.22%  "
Set-TextValue $ws.Range("D18") "11.49"
Set-TextValue $ws.Range("E18") "  +1.19%  "
Set-TextValue $ws.Range("E19") "  +3.92%  "
Set-TextValue $ws.Range("D20") "342.49"
Set-TextValue $ws.Range("E20") "  +0.44%  "
Set-TextValue $ws.Range("E21") "  +4.42%  "
Set-TextValue $ws.Range("E22") "  -0.13%  "
Set-TextValue $ws.Range("D23") "5.60"
Set-TextValue $ws.Range("E23") "  -2.77%  "
Set-TextValue $ws.Range("D24") "66.85"
Set-TextValue $ws.Range("E24") "  -1.11%  "
Set-TextValue $ws.Range("D25") "1.69"
Set-TextValue $ws.Range("E25") "  +0.89%  "
Set-TextValue $ws.Range("D26") "9.03"
Set-TextValue $ws.Range("E26") "  +7.43%  "
Set-TextValue $ws.Range("D27") "1.53"
Set-TextValue $ws.Range("E27") "  -1.39%  "
Set-TextValue $ws.Range("D28") "0.165"
Set-TextValue $ws.Range("E28") "  -0.54%  "
Set-TextValue $ws.Range("D29") "549.31"
Set-TextValue $ws.Range("E29") "  -0.65%  "
Set-TextValue $ws.Range("D30") "0.998"
Set-TextValue $ws.Range("E30") "  -0.12%  "
Set-TextValue $ws.Range("E31") "  +1.03%  "
Set-TextValue $ws.Range("D32") "2.05"
Set-TextValue $ws.Range("E32") "  +3.66%  "
Set-TextValue $ws.Range("D33") "1.76"
Set-TextValue $ws.Range("E33") "  -3.12%  "
Set-TextValue $ws.Range("D34") "0.0₃0818"
Set-TextValue $ws.Range("E34") "  +1.53%  "
Set-TextValue $ws.Range("D35") "5.19"
Set-TextValue $ws.Range("E35") "  +6.45%  "
Set-TextValue $ws.Range("D36") "167.27"
Set-TextValue $ws.Range("E36") "  -4.60%  "
Set-TextValue $ws.Range("D37") "0.407"
Set-TextValue $ws.Range("E37") "  +1.39%  "
Set-TextValue $ws.Range("E38") "  -0.03%  "
Set-TextValue $ws.Range("B39") "EthereumClassic"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D39") "19.14"
Set-TextValue $ws.Range("E39") "  +0.45%  "
Set-TextValue $ws.Range("B40") "Stacks"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D40") "1.91"
Set-TextValue $ws.Range("E40") "  +7.36%  "
Set-TextValue $ws.Range("E41") "  +0.02%  "
Set-TextValue $ws.Range("D42") "169.35"
Set-TextValue $ws.Range("E42") "  -0.49%  "
Set-TextValue $ws.Range("E43") "  +1.52%  "
Set-TextValue $ws.Range("D44") "22.69"
Set-TextValue $ws.Range("E44") "  +1.61%  "
Set-TextValue $ws.Range("D45") "0.0577"
Set-TextValue $ws.Range("E45") "  +4.91%  "
Set-TextValue $ws.Range("D46") "0.630"
Set-TextValue $ws.Range("E46") "  +0.36%  "
Set-TextValue $ws.Range("D47") "0.0249"
Set-TextValue $ws.Range("E47") "  +4.51%  "
Set-TextValue $ws.Range("D48") "0.0964"
Set-TextValue $ws.Range("E48") "  +0.56%  "
Set-TextValue $ws.Range("E49") "  +1.44%  "
Set-TextValue $ws.Range("D50") "1.88"
Set-TextValue $ws.Range("E50") "  +10.49%  "
Set-TextValue $ws.Range("E51") "  -0.56%  "
